# Append the three new "Mobile" command rows (LanguageKey / en_us pairs)
# to the bottom of the localization table on Sheet1, directly below the
# existing last row (row 21: "$excuteCustomMethod" / "Execute Custom Method").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("`$selectValueInMobileSelect", "Select value in Mobile Select DropDown"),
    @("`$selectMobileRadioButton",   "Select mobile radio Button"),
    @("`$setDateInMobile",           "Set date in Mobile")
)

$startRow = 22
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}

# Scroll the viewport down a bit (matches the author having scrolled to
# row 3 before saving) and leave the selection on the first empty cell
# below the newly-entered data, just like Excel does after typing values
# down a column and pressing Enter.
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 3
    $win.ScrollColumn = 1
} catch {
    # Non-fatal: scroll position is cosmetic only.
}

$ws.Range("B25").Select()
